# [#9079] Set up PWA/SW architecture and migrate to Jest for testing (#9263)
#
# The deck's "testing stack" diagram has a text box ("TextBox 75") that lists
# the testing tools used, one per paragraph ("TestNG," / "Karma/Jasmine").
# The commit replaces the old JS test runner/framework with Jest, so the
# "Karma/Jasmine" paragraph becomes "Jest".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Find the shape by name rather than a hard-coded index so the script keeps
# working even if shape ordering on the slide ever changes.
$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "TextBox 75") {
        $target = $shp
    }
}

$textRange = $target.TextFrame.TextRange

# Locate the paragraph that reads "Karma/Jasmine" and rewrite just that run's
# text to "Jest", leaving every other paragraph/run (and their formatting)
# untouched.
$paragraphCount = $textRange.Paragraphs().Count
for ($j = 1; $j -le $paragraphCount; $j++) {
    $paragraph = $textRange.Paragraphs($j)
    if ($paragraph.Text -eq "Karma/Jasmine") {
        $paragraph.Runs(1).Text = "Jest"
    }
}
